$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.0543
$ws.Range("B9").Value = 8.659400000000003
$ws.Range("C11").Value = -13.2061
$ws.Range("B18").Value = 4.599100000000004
$ws.Range("B20").Value = 5.593800000000001
